$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.1282389760017395
$ws.Cells.Item(2, 2).Value = 0.9628564715385437
$ws.Cells.Item(2, 3).Value = 0.1785556823015213
$ws.Cells.Item(2, 4).Value = 0.9517366886138916
$ws.Cells.Item(3, 1).Value = 0.03223736211657524
$ws.Cells.Item(3, 2).Value = 0.992206871509552
$ws.Cells.Item(3, 3).Value = 0.2270778119564056
$ws.Cells.Item(3, 4).Value = 0.9411821365356445
$ws.Cells.Item(4, 1).Value = 0.02203092165291309
$ws.Cells.Item(4, 2).Value = 0.9930191040039062
$ws.Cells.Item(4, 3).Value = 0.1386857032775879
$ws.Cells.Item(4, 4).Value = 0.9526962041854858
$ws.Cells.Item(5, 1).Value = 0.01683168485760689
$ws.Cells.Item(5, 2).Value = 0.9940289258956909
$ws.Cells.Item(5, 3).Value = 0.07311633229255676
$ws.Cells.Item(5, 4).Value = 0.9596046805381775
$ws.Cells.Item(6, 1).Value = 0.01440321281552315
$ws.Cells.Item(6, 2).Value = 0.9941606521606445
$ws.Cells.Item(6, 3).Value = 0.1242493987083435
$ws.Cells.Item(6, 4).Value = 0.9579735398292542
$ws.Cells.Item(7, 1).Value = 0.01380656287074089
$ws.Cells.Item(7, 2).Value = 0.9937655329704285
$ws.Cells.Item(7, 3).Value = 0.121691606938839
$ws.Cells.Item(7, 4).Value = 0.954327404499054
$ws.Cells.Item(8, 1).Value = 0.01195869967341423
$ws.Cells.Item(8, 2).Value = 0.9963339567184448
$ws.Cells.Item(8, 3).Value = 0.09574044495820999
$ws.Cells.Item(8, 4).Value = 0.957685649394989
$ws.Cells.Item(9, 1).Value = 0.01050076633691788
$ws.Cells.Item(9, 2).Value = 0.9967730045318604
$ws.Cells.Item(9, 3).Value = 0.07743805646896362
$ws.Cells.Item(9, 4).Value = 0.9625791311264038
$ws.Cells.Item(10, 1).Value = 0.01018921751528978
$ws.Cells.Item(10, 2).Value = 0.9965754151344299
$ws.Cells.Item(10, 3).Value = 0.1396393626928329
$ws.Cells.Item(10, 4).Value = 0.9459796547889709
$ws.Cells.Item(11, 1).Value = 0.009051336906850338
$ws.Cells.Item(11, 2).Value = 0.996816873550415
$ws.Cells.Item(11, 3).Value = 0.08929783850908279
$ws.Cells.Item(11, 4).Value = 0.9665131568908691
$ws.Cells.Item(12, 1).Value = 0.009744052775204182
$ws.Cells.Item(12, 2).Value = 0.9965095520019531
$ws.Cells.Item(12, 3).Value = 0.08521924167871475
$ws.Cells.Item(12, 4).Value = 0.9632508158683777
$ws.Cells.Item(13, 1).Value = 0.008984826505184174
$ws.Cells.Item(13, 2).Value = 0.9966632127761841
$ws.Cells.Item(13, 3).Value = 0.08180578798055649
$ws.Cells.Item(13, 4).Value = 0.966225266456604
$ws.Cells.Item(14, 1).Value = 0.009071167558431625
$ws.Cells.Item(14, 2).Value = 0.9966193437576294
$ws.Cells.Item(14, 3).Value = 0.09515927731990814
$ws.Cells.Item(14, 4).Value = 0.9754365682601929
$ws.Cells.Item(15, 1).Value = 0.008991554379463196
$ws.Cells.Item(15, 2).Value = 0.9965754151344299
$ws.Cells.Item(15, 3).Value = 0.08974380791187286
$ws.Cells.Item(15, 4).Value = 0.9614277482032776
$ws.Cells.Item(16, 1).Value = 0.008799412287771702
$ws.Cells.Item(16, 2).Value = 0.9968827366828918
$ws.Cells.Item(16, 3).Value = 0.1180259585380554
$ws.Cells.Item(16, 4).Value = 0.9680483341217041
$ws.Cells.Item(17, 1).Value = 0.008448570035398006
$ws.Cells.Item(17, 2).Value = 0.9967071413993835
$ws.Cells.Item(17, 3).Value = 0.1644821614027023
$ws.Cells.Item(17, 4).Value = 0.9670888781547546
$ws.Cells.Item(18, 1).Value = 0.009483812376856804
$ws.Cells.Item(18, 2).Value = 0.9966412782669067
$ws.Cells.Item(18, 3).Value = 0.05399708077311516
$ws.Cells.Item(18, 4).Value = 0.9668009877204895
$ws.Cells.Item(19, 1).Value = 0.008651969023048878
$ws.Cells.Item(19, 2).Value = 0.9965973496437073
$ws.Cells.Item(19, 3).Value = 0.05299533531069756
$ws.Cells.Item(19, 4).Value = 0.9697754979133606
$ws.Cells.Item(20, 1).Value = 0.008957092650234699
$ws.Cells.Item(20, 2).Value = 0.9967510104179382
$ws.Cells.Item(20, 3).Value = 0.1048153787851334
$ws.Cells.Item(20, 4).Value = 0.967568576335907
$ws.Cells.Item(21, 1).Value = 0.008600062690675259
$ws.Cells.Item(21, 2).Value = 0.9966193437576294
$ws.Cells.Item(21, 3).Value = 0.05450410023331642
$ws.Cells.Item(21, 4).Value = 0.9722701907157898
$ws.Cells.Item(22, 1).Value = 0.009247648529708385
$ws.Cells.Item(22, 2).Value = 0.9965754151344299
$ws.Cells.Item(22, 3).Value = 0.04101405665278435
$ws.Cells.Item(22, 4).Value = 0.9781231880187988
$ws.Cells.Item(23, 1).Value = 0.008323425427079201
$ws.Cells.Item(23, 2).Value = 0.9967949390411377
$ws.Cells.Item(23, 3).Value = 0.03869495540857315
$ws.Cells.Item(23, 4).Value = 0.9786029458045959
$ws.Cells.Item(24, 1).Value = 0.009087215177714825
$ws.Cells.Item(24, 2).Value = 0.9965095520019531
$ws.Cells.Item(24, 3).Value = 0.04440681263804436
$ws.Cells.Item(24, 4).Value = 0.978411078453064
$ws.Cells.Item(25, 1).Value = 0.009995067492127419
$ws.Cells.Item(25, 2).Value = 0.9963119626045227
$ws.Cells.Item(25, 3).Value = 0.09454083442687988
$ws.Cells.Item(25, 4).Value = 0.9694876074790955
$ws.Cells.Item(26, 1).Value = 0.008116000331938267
$ws.Cells.Item(26, 2).Value = 0.996904730796814
$ws.Cells.Item(26, 3).Value = 0.1313388347625732
$ws.Cells.Item(26, 4).Value = 0.9692957401275635
$ws.Cells.Item(27, 1).Value = 0.00849370751529932
$ws.Cells.Item(27, 2).Value = 0.9967510104179382
$ws.Cells.Item(27, 3).Value = 0.1361873596906662
$ws.Cells.Item(27, 4).Value = 0.9710228443145752
$ws.Cells.Item(28, 1).Value = 0.008890886791050434
$ws.Cells.Item(28, 2).Value = 0.996421754360199
$ws.Cells.Item(28, 3).Value = 0.1485448777675629
$ws.Cells.Item(28, 4).Value = 0.9641143679618835
$ws.Cells.Item(29, 1).Value = 0.008416827768087387
$ws.Cells.Item(29, 2).Value = 0.9967949390411377
$ws.Cells.Item(29, 3).Value = 0.02545074373483658
$ws.Cells.Item(29, 4).Value = 0.9865669012069702
$ws.Cells.Item(30, 1).Value = 0.008944380097091198
$ws.Cells.Item(30, 2).Value = 0.9963339567184448
$ws.Cells.Item(30, 3).Value = 0.01440200302749872
$ws.Cells.Item(30, 4).Value = 0.9949145913124084
$ws.Cells.Item(31, 1).Value = 0.009300052188336849
$ws.Cells.Item(31, 2).Value = 0.9961802959442139
$ws.Cells.Item(31, 3).Value = 0.01399672497063875
$ws.Cells.Item(31, 4).Value = 0.9948186278343201
$ws.Cells.Item(32, 1).Value = 0.009228548035025597
$ws.Cells.Item(32, 2).Value = 0.9970144629478455
$ws.Cells.Item(32, 3).Value = 0.01115142926573753
$ws.Cells.Item(32, 4).Value = 0.9965457916259766
$ws.Cells.Item(33, 1).Value = 0.007855038158595562
$ws.Cells.Item(33, 2).Value = 0.9968827366828918
$ws.Cells.Item(33, 3).Value = 0.009478968568146229
$ws.Cells.Item(33, 4).Value = 0.997313380241394
$ws.Cells.Item(34, 1).Value = 0.008763711899518967
$ws.Cells.Item(34, 2).Value = 0.9966193437576294
$ws.Cells.Item(34, 3).Value = 0.01785976625978947
$ws.Cells.Item(34, 4).Value = 0.99299556016922
$ws.Cells.Item(35, 1).Value = 0.008519371971487999
$ws.Cells.Item(35, 2).Value = 0.9965973496437073
$ws.Cells.Item(35, 3).Value = 0.01033820491284132
$ws.Cells.Item(35, 4).Value = 0.9979850053787231
$ws.Cells.Item(36, 1).Value = 0.008793055079877377
$ws.Cells.Item(36, 2).Value = 0.9963998198509216
$ws.Cells.Item(36, 3).Value = 0.009882328100502491
$ws.Cells.Item(36, 4).Value = 0.9977931380271912
$ws.Cells.Item(37, 1).Value = 0.008644300512969494
$ws.Cells.Item(37, 2).Value = 0.9967510104179382
$ws.Cells.Item(37, 3).Value = 0.01758559979498386
$ws.Cells.Item(37, 4).Value = 0.9953943490982056
$ws.Cells.Item(38, 1).Value = 0.009450608864426613
$ws.Cells.Item(38, 2).Value = 0.9962022304534912
$ws.Cells.Item(38, 3).Value = 0.01575752533972263
$ws.Cells.Item(38, 4).Value = 0.9969295859336853
$ws.Cells.Item(39, 1).Value = 0.008584649302065372
$ws.Cells.Item(39, 2).Value = 0.9965754151344299
$ws.Cells.Item(39, 3).Value = 0.02285070158541203
$ws.Cells.Item(39, 4).Value = 0.9944348335266113
$ws.Cells.Item(40, 1).Value = 0.008455325849354267
$ws.Cells.Item(40, 2).Value = 0.9966412782669067
$ws.Cells.Item(40, 3).Value = 0.006157251540571451
$ws.Cells.Item(40, 4).Value = 0.9989445209503174
$ws.Cells.Item(41, 1).Value = 0.007854791358113289
$ws.Cells.Item(41, 2).Value = 0.9970363974571228
$ws.Cells.Item(41, 3).Value = 0.00817712489515543
$ws.Cells.Item(41, 4).Value = 0.9982728958129883
$ws.Cells.Item(42, 1).Value = 0.008207334205508232
$ws.Cells.Item(42, 2).Value = 0.9967290759086609
$ws.Cells.Item(42, 3).Value = 0.01105641294270754
$ws.Cells.Item(42, 4).Value = 0.9980809688568115
$ws.Cells.Item(43, 1).Value = 0.008842155337333679
$ws.Cells.Item(43, 2).Value = 0.9964436888694763
$ws.Cells.Item(43, 3).Value = 0.01407072227448225
$ws.Cells.Item(43, 4).Value = 0.9957781434059143
$ws.Cells.Item(44, 1).Value = 0.008421175181865692
$ws.Cells.Item(44, 2).Value = 0.9965973496437073
$ws.Cells.Item(44, 3).Value = 0.008083195425570011
$ws.Cells.Item(44, 4).Value = 0.998656690120697
$ws.Cells.Item(45, 1).Value = 0.007668427191674709
$ws.Cells.Item(45, 2).Value = 0.9970144629478455
$ws.Cells.Item(45, 3).Value = 0.007804097607731819
$ws.Cells.Item(45, 4).Value = 0.9984647631645203
$ws.Cells.Item(46, 1).Value = 0.009726038202643394
$ws.Cells.Item(46, 2).Value = 0.9967949390411377
$ws.Cells.Item(46, 3).Value = 0.006464036181569099
$ws.Cells.Item(46, 4).Value = 0.9988486170768738
$ws.Cells.Item(47, 1).Value = 0.00785375852137804
$ws.Cells.Item(47, 2).Value = 0.9968388676643372
$ws.Cells.Item(47, 3).Value = 0.007867932319641113
$ws.Cells.Item(47, 4).Value = 0.9988486170768738
$ws.Cells.Item(48, 1).Value = 0.007840006612241268
$ws.Cells.Item(48, 2).Value = 0.9968388676643372
$ws.Cells.Item(48, 3).Value = 0.007646213285624981
$ws.Cells.Item(48, 4).Value = 0.9985607266426086
$ws.Cells.Item(49, 1).Value = 0.008642685599625111
$ws.Cells.Item(49, 2).Value = 0.9966193437576294
$ws.Cells.Item(49, 3).Value = 0.007394760381430387
$ws.Cells.Item(49, 4).Value = 0.998656690120697
$ws.Cells.Item(50, 1).Value = 0.008498113602399826
$ws.Cells.Item(50, 2).Value = 0.9965314865112305
$ws.Cells.Item(50, 3).Value = 0.007553883362561464
$ws.Cells.Item(50, 4).Value = 0.9984647631645203
$ws.Cells.Item(51, 1).Value = 0.008400174789130688
$ws.Cells.Item(51, 2).Value = 0.9965973496437073
$ws.Cells.Item(51, 3).Value = 0.01532774604856968
$ws.Cells.Item(51, 4).Value = 0.9945307970046997
